# UserDashboard seed workbook update:
#   "v1/welcome call now returns UserDetail" -> the dashboard list row that
#   used to hold "Danny Business Money" is renamed to "Test User Dash" and
#   moved to the top of the list (directly under Tom's / Bethany's / Danny's
#   personal boards, which shift down one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: now the renamed "Test User Dash" entry
$ws.Range("B2").Value = "Test User Dash"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1

# Row 3: was row 2 (Tom's Personal Money Tracker Board)
$ws.Range("B3").Value = "Tom's Personal Money Tracker Board"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 11

# Row 4: was row 3 (Bethany's Bread & Buns Budget Board)
$ws.Range("B4").Value = "Bethany's Bread & Buns Budget Board"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 12

# Row 5: was row 4 (Danny Personal Money)
$ws.Range("B5").Value = "Danny Personal Money"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 13

# Column B widened to fit the longest label now in the list
$ws.Columns("B").ColumnWidth = 31.4

# Last selected cell moved to F7
$ws.Range("F7").Select()
